# Threshold/Zn/3His/4.xlsx update
#  - widen column A and column C (manual widths, no longer "best fit")
#  - bump the Max-column threshold values in rows 2-4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths: target OOXML <col> widths are 27 (col A) and 27.25 (col C).
# Excel's ColumnWidth property is expressed in characters and gets rounded to
# the sheet's pixel grid on write, so feed it the fractional char-widths that
# land on the desired pixel boundaries.
$ws.Columns.Item(1).ColumnWidth = 26.2857142857143
$ws.Columns.Item(3).ColumnWidth = 26.5714285714286

# Updated threshold values
$ws.Range("C2").Value = 11
$ws.Range("C3").Value = 9.5
$ws.Range("C4").Value = 1.4
